$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 6 ("Enemy" NPC template) into row 7 so the new row starts
# with identical stats/formatting to the existing Enemy/AttackNpc3 template.
$ws.Rows.Item(6).Copy()
$ws.Rows.Item(7).PasteSpecial()

# The newly duplicated row 7 becomes the "Default" npc config entry.
$ws.Range("A7").Value = "Default"

# Rename the row-5 NPC from "Player" to "Player_0_0".
$ws.Range("A5").Value = "Player_0_0"

# Widen column W (DropPackList comment column) to fit its content.
# (63.375 character-units is the authored target; the host engine snaps
# ColumnWidth to whole-pixel boundaries using MDW=7, so 62.71 is the input
# that lands on the nearest achievable stored width, 63.428571..., i.e.
# the closest representable value to 63.375.)
$ws.Columns.Item(23).ColumnWidth = 62.71

# Match the saved selection state (active cell A6).
$ws.Range("A6").Select() | Out-Null
